{"js": "// Closes #29 - apply strikethrough formatting to the completed to-do item\n// \"Code review both backend and frontend.\" (including the paragraph mark),\n// mirroring a user selecting the whole list item and toggling Strikethrough.\n\nconst searchResults = context.document.body.search(\n  \"Code review both backend and frontend.\",\n  { matchCase: true, matchWholeWord: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Paragraph \"Code review both backend and frontend.\" was not found.');\n}\n\n// Grab the paragraph that contains the matched text so the paragraph mark\n// (the trailing pilcrow) is included in the formatting change, not just the\n// run text.\nconst targetParagraph = searchResults.items[0].paragraphs.getFirst();\n\n// Applying strikethrough at the paragraph level (via its Font) marks both\n// the paragraph's own run(s) and the paragraph mark's run properties,\n// exactly like selecting the full line (incl. end-of-paragraph) in the UI\n// and pressing the Strikethrough toggle.\ntargetParagraph.font.strikeThrough = true;\n\nawait context.sync();\n", "ps1": "# Closes #29 - apply strikethrough formatting to the completed to-do item\n# \"Code review both backend and frontend.\" (including the paragraph mark),\n# mirroring a user selecting the whole list item and toggling Strikethrough.\n\n$d = $word.ActiveDocument\n\n$searchText = \"Code review both backend and frontend.\"\n\n$finder = $d.Content\n$found = $finder.Find.Execute($searchText)\nif (-not $found) {\n    throw \"Paragraph '$searchText' was not found.\"\n}\n\n$foundStart = $finder.Start\n$foundEnd = $finder.End\n\n# Resolve the Paragraph object that contains the found text so that\n# formatting the Paragraph's Range also stamps the paragraph mark (pilcrow)\n# run properties, not just the visible run text.\n$targetParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $foundStart -and $p.Range.End -ge $foundEnd) {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($null -eq $targetParagraph) {\n    throw \"Could not resolve paragraph for '$searchText'.\"\n}\n\n$targetParagraph.Range.Font.StrikeThrough = 1\n"}
